$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.173.92"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.895.29"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.80"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5370"
$ws.Range("E7").Value = "  +3.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3786"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07261"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.87"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8947"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08166"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.47"
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.330"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.819.41"
$ws.Range("E15").Value = "  -6.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.82"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008625"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.060.69"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.024"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.452"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.67"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.281"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.28"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.740"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.01"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.806"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.661"
$ws.Range("E30").Value = "  -4.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09158"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8178"
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05026"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.211"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.015"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.282"
$ws.Range("E36").Value = "  -4.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.658"
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5939"
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01980"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.074"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.221"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.606"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.96"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5059"
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1520"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.15"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.622"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.83"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06066"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.70"
$ws.Range("E51").Value = "  -1.72%  "
